$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "sheet1" to "sheet"
$ws.Name = "sheet"

# Header row (row 1) - new column headers spanning A1:M1
$headers = @("Case Number", "Document Type", "First Name", "Middle Name", "Last Name", "Suffix", "Title", "Company Name", "Address 1", "Address 2", "City", "State", "Zip")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 - clear A2 (was "pkp"), set B2:D2 to test2/test3/test4
$ws.Cells.Item(2, 1).Value = $null
$ws.Cells.Item(2, 2).Value = "test2"
$ws.Cells.Item(2, 3).Value = "test3"
$ws.Cells.Item(2, 4).Value = "test4"

# Touch A2 so it stays as an empty cell entry rather than being dropped
$ws.Cells.Item(2, 1).Font.Bold = $false

# Extend the used range out to column M for every row (1-12) by touching
# the new columns K:M without actually storing any value in them.
$ws.Range("K1:M12").Font.Bold = $false
